$wb = $excel.ActiveWorkbook

# The workbook stores the same event list on both the "展览" sheet and the
# "全部类型" aggregate sheet. Row 7 (CCAC) and row 8 (AEO) need their
# "想去人数" (column F) counts bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F7").Value = 34
    $ws.Range("F8").Value = 181
}
